$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 16.75  # H6: 22.75 -> 16.75
$ws.Cells.Item(6, 9).Value = 16.75  # I6: 22.75 -> 16.75
$ws.Cells.Item(6, 11).Value = 50.25  # K6: 68.25 -> 50.25
$ws.Cells.Item(6, 13).Value = 61.75  # M6: 43.75 -> 61.75

$ws.Cells.Item(17, 8).Value = 1309397.2  # H17: 1627291.2 -> 1309397.2
$ws.Cells.Item(17, 10).Value = 1309397.2  # J17: 1627291.2 -> 1309397.2
$ws.Cells.Item(17, 12).Value = 3928191.6  # L17: 4881873.6 -> 3928191.6
$ws.Cells.Item(17, 14).Value = -3928527.6  # N17: -4882209.6 -> -3928527.6

$ws.Cells.Item(28, 8).Value = 428.75  # H28: 423.85715 -> 428.75
$ws.Cells.Item(28, 9).Value = 429.60715  # I28: 423.93103 -> 429.60715
$ws.Cells.Item(28, 10).Value = 422.75  # J28: 423.5 -> 422.75
$ws.Cells.Item(28, 11).Value = 429.60715  # K28: 423.93103 -> 429.60715
$ws.Cells.Item(28, 12).Value = 422.75  # L28: 423.5 -> 422.75
$ws.Cells.Item(28, 13).Value = 55.39285000000001  # M28: 61.06896999999998 -> 55.39285000000001
$ws.Cells.Item(28, 14).Value = -1392.75  # N28: -1393.5 -> -1392.75

$ws.Cells.Item(29, 8).Value = 893.7059  # H29: 1519 -> 893.7059
$ws.Cells.Item(29, 9).Value = 707.1539  # I29: 1519 -> 707.1539
$ws.Cells.Item(29, 10).Value = 1500  # J29: 0 -> 1500
$ws.Cells.Item(29, 11).Value = 2121.4617  # K29: 4557 -> 2121.4617
$ws.Cells.Item(29, 12).Value = 4500  # L29: 0 -> 4500
$ws.Cells.Item(29, 13).Value = -1840.4617  # M29: -4276 -> -1840.4617
$ws.Cells.Item(29, 14).Value = -5062  # N29: (new cell) -> -5062

$ws.Cells.Item(38, 8).Value = 1137.5834  # H38: 1230.8182 -> 1137.5834
$ws.Cells.Item(38, 9).Value = 93.833336  # I38: 101.63636 -> 93.833336
$ws.Cells.Item(38, 10).Value = 2181.3333  # J38: 2360 -> 2181.3333
$ws.Cells.Item(38, 11).Value = 281.500008  # K38: 304.90908 -> 281.500008
$ws.Cells.Item(38, 12).Value = 6543.999899999999  # L38: 7080 -> 6543.999899999999
$ws.Cells.Item(38, 13).Value = 90.49999200000002  # M38: 67.09091999999998 -> 90.49999200000002
$ws.Cells.Item(38, 14).Value = -7287.999899999999  # N38: -7824 -> -7287.999899999999

$ws.Cells.Item(98, 8).Value = 2251.6191  # H98: 2361.75 -> 2251.6191
$ws.Cells.Item(98, 9).Value = 2314.2  # I98: 2433.4211 -> 2314.2
$ws.Cells.Item(98, 11).Value = 2314.2  # K98: 2433.4211 -> 2314.2
$ws.Cells.Item(98, 13).Value = -816.1999999999998  # M98: -935.4211 -> -816.1999999999998

$ws.Cells.Item(107, 8).Value = 659.1818  # H107: 684.6 -> 659.1818
$ws.Cells.Item(107, 9).Value = 601.6667  # I107: 568.3333 -> 601.6667
$ws.Cells.Item(107, 10).Value = 728.2  # J107: 859 -> 728.2
$ws.Cells.Item(107, 11).Value = 601.6667  # K107: 568.3333 -> 601.6667
$ws.Cells.Item(107, 12).Value = 728.2  # L107: 859 -> 728.2
$ws.Cells.Item(107, 13).Value = 1318.3333  # M107: 1351.6667 -> 1318.3333
$ws.Cells.Item(107, 14).Value = -4568.2  # N107: -4699 -> -4568.2

$ws.Cells.Item(122, 8).Value = 2251.6191  # H122: 2361.75 -> 2251.6191
$ws.Cells.Item(122, 9).Value = 2314.2  # I122: 2433.4211 -> 2314.2
$ws.Cells.Item(122, 11).Value = 6942.599999999999  # K122: 7300.263300000001 -> 6942.599999999999
$ws.Cells.Item(122, 13).Value = -4492.599999999999  # M122: -4850.263300000001 -> -4492.599999999999

$ws.Cells.Item(131, 8).Value = 2829.111  # H131: 2519.625 -> 2829.111
$ws.Cells.Item(131, 9).Value = 560.3333  # I131: 693.6667 -> 560.3333
$ws.Cells.Item(131, 10).Value = 7366.6665  # J131: 7997.5 -> 7366.6665
$ws.Cells.Item(131, 11).Value = 1680.9999  # K131: 2081.0001 -> 1680.9999
$ws.Cells.Item(131, 12).Value = 22099.9995  # L131: 23992.5 -> 22099.9995
$ws.Cells.Item(131, 13).Value = 3359.0001  # M131: 2958.9999 -> 3359.0001
$ws.Cells.Item(131, 14).Value = -32179.9995  # N131: -34072.5 -> -32179.9995

$ws.Cells.Item(132, 8).Value = 10558.333  # H132: 10778.906 -> 10558.333
$ws.Cells.Item(132, 10).Value = 28376.25  # J132: 36668.332 -> 28376.25
$ws.Cells.Item(132, 12).Value = 85128.75  # L132: 110004.996 -> 85128.75
$ws.Cells.Item(132, 14).Value = -90188.75  # N132: -115064.996 -> -90188.75

$ws.Cells.Item(135, 8).Value = 1508.8  # H135: 1563.4584 -> 1508.8
$ws.Cells.Item(135, 9).Value = 1769.5  # I135: 1852.2632 -> 1769.5
$ws.Cells.Item(135, 11).Value = 15925.5  # K135: 16670.3688 -> 15925.5
$ws.Cells.Item(135, 13).Value = -13390.5  # M135: -14135.3688 -> -13390.5

$ws.Cells.Item(137, 8).Value = 7081.1606  # H137: 7214.309 -> 7081.1606
$ws.Cells.Item(137, 9).Value = 1150.862  # I137: 1166.1034 -> 1150.862
$ws.Cells.Item(137, 10).Value = 13450.741  # J137: 13960.385 -> 13450.741
$ws.Cells.Item(137, 11).Value = 3452.586  # K137: 3498.3102 -> 3452.586
$ws.Cells.Item(137, 12).Value = 40352.223  # L137: 41881.155 -> 40352.223
$ws.Cells.Item(137, 13).Value = -902.5860000000002  # M137: -948.3101999999999 -> -902.5860000000002
$ws.Cells.Item(137, 14).Value = -45452.223  # N137: -46981.155 -> -45452.223

$ws.Cells.Item(138, 8).Value = 2854.875  # H138: 2900.875 -> 2854.875
$ws.Cells.Item(138, 9).Value = 3764.5557  # I138: 4122.625 -> 3764.5557
$ws.Cells.Item(138, 10).Value = 2309.0667  # J138: 2290 -> 2309.0667
$ws.Cells.Item(138, 11).Value = 11293.6671  # K138: 12367.875 -> 11293.6671
$ws.Cells.Item(138, 12).Value = 6927.2001  # L138: 6870 -> 6927.2001
$ws.Cells.Item(138, 13).Value = -6153.667099999999  # M138: -7227.875 -> -6153.667099999999
$ws.Cells.Item(138, 14).Value = -17207.2001  # N138: -17150 -> -17207.2001

$ws.Cells.Item(141, 8).Value = 2943.182  # H141: 3509.111 -> 2943.182
$ws.Cells.Item(141, 9).Value = 1459.75  # I141: 1814.1666 -> 1459.75
$ws.Cells.Item(141, 11).Value = 4379.25  # K141: 5442.4998 -> 4379.25
$ws.Cells.Item(141, 13).Value = 800.75  # M141: -262.4997999999996 -> 800.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1508.3572  # H5: 1421.6 -> 1508.3572
$ws.Cells.Item(5, 9).Value = 2481.625  # I5: 2216.5557 -> 2481.625
$ws.Cells.Item(5, 10).Value = 210.66667  # J5: 229.16667 -> 210.66667
$ws.Cells.Item(5, 11).Value = 2481.625  # K5: 2216.5557 -> 2481.625
$ws.Cells.Item(5, 12).Value = 210.66667  # L5: 229.16667 -> 210.66667
$ws.Cells.Item(5, 13).Value = -2369.625  # M5: -2104.5557 -> -2369.625
$ws.Cells.Item(5, 14).Value = -434.66667  # N5: -453.16667 -> -434.66667

$ws.Cells.Item(45, 8).Value = 2578.4736  # H45: 2616.2222 -> 2578.4736
$ws.Cells.Item(45, 10).Value = 2859.125  # J45: 2996.2856 -> 2859.125
$ws.Cells.Item(45, 12).Value = 2859.125  # L45: 2996.2856 -> 2859.125
$ws.Cells.Item(45, 14).Value = -3613.125  # N45: -3750.2856 -> -3613.125

$ws.Cells.Item(51, 8).Value = 0  # H51: 39900 -> 0
$ws.Cells.Item(51, 10).Value = 0  # J51: 39900 -> 0
$ws.Cells.Item(51, 12).Value = 0  # L51: 39900 -> 0
$ws.Cells.Item(51, 14).ClearContents() | Out-Null  # N51: remove (was -41412)

$ws.Cells.Item(61, 8).Value = 9437.208000000001  # H61: 10577.617 -> 9437.208000000001
$ws.Cells.Item(61, 9).Value = 2005.6842  # I61: 2637.2307 -> 2005.6842
$ws.Cells.Item(61, 10).Value = 13590.117  # J61: 13613.647 -> 13590.117
$ws.Cells.Item(61, 11).Value = 2005.6842  # K61: 2637.2307 -> 2005.6842
$ws.Cells.Item(61, 12).Value = 13590.117  # L61: 13613.647 -> 13590.117
$ws.Cells.Item(61, 13).Value = -1793.6842  # M61: -2425.2307 -> -1793.6842
$ws.Cells.Item(61, 14).Value = -14014.117  # N61: -14037.647 -> -14014.117

$ws.Cells.Item(74, 8).Value = 30668.5  # H74: 28816.533 -> 30668.5
$ws.Cells.Item(74, 9).Value = 1822  # I74: 1964.4 -> 1822
$ws.Cells.Item(74, 10).Value = 52303.375  # J74: 42242.6 -> 52303.375
$ws.Cells.Item(74, 11).Value = 1822  # K74: 1964.4 -> 1822
$ws.Cells.Item(74, 12).Value = 52303.375  # L74: 42242.6 -> 52303.375
$ws.Cells.Item(74, 13).Value = -948  # M74: -1090.4 -> -948
$ws.Cells.Item(74, 14).Value = -54051.375  # N74: -43990.6 -> -54051.375

$ws.Cells.Item(77, 8).Value = 30668.5  # H77: 28816.533 -> 30668.5
$ws.Cells.Item(77, 9).Value = 1822  # I77: 1964.4 -> 1822
$ws.Cells.Item(77, 10).Value = 52303.375  # J77: 42242.6 -> 52303.375
$ws.Cells.Item(77, 11).Value = 9110  # K77: 9822 -> 9110
$ws.Cells.Item(77, 12).Value = 261516.875  # L77: 211213 -> 261516.875
$ws.Cells.Item(77, 13).Value = -4742  # M77: -5454 -> -4742
$ws.Cells.Item(77, 14).Value = -270252.875  # N77: -219949 -> -270252.875

$ws.Cells.Item(136, 8).Value = 9437.208000000001  # H136: 10577.617 -> 9437.208000000001
$ws.Cells.Item(136, 9).Value = 2005.6842  # I136: 2637.2307 -> 2005.6842
$ws.Cells.Item(136, 10).Value = 13590.117  # J136: 13613.647 -> 13590.117
$ws.Cells.Item(136, 11).Value = 6017.0526  # K136: 7911.6921 -> 6017.0526
$ws.Cells.Item(136, 12).Value = 40770.351  # L136: 40840.94100000001 -> 40770.351
$ws.Cells.Item(136, 13).Value = -3467.0526  # M136: -5361.6921 -> -3467.0526
$ws.Cells.Item(136, 14).Value = -45870.351  # N136: -45940.94100000001 -> -45870.351

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1508.3572  # H4: 1421.6 -> 1508.3572
$ws.Cells.Item(4, 9).Value = 2481.625  # I4: 2216.5557 -> 2481.625
$ws.Cells.Item(4, 10).Value = 210.66667  # J4: 229.16667 -> 210.66667
$ws.Cells.Item(4, 11).Value = 2481.625  # K4: 2216.5557 -> 2481.625
$ws.Cells.Item(4, 12).Value = 210.66667  # L4: 229.16667 -> 210.66667
$ws.Cells.Item(4, 13).Value = -2366.625  # M4: -2101.5557 -> -2366.625
$ws.Cells.Item(4, 14).Value = -440.66667  # N4: -459.16667 -> -440.66667

$ws.Cells.Item(110, 8).Value = 100000  # H110: 67275 -> 100000
$ws.Cells.Item(110, 10).Value = 100000  # J110: 67275 -> 100000
$ws.Cells.Item(110, 12).Value = 100000  # L110: 67275 -> 100000
$ws.Cells.Item(110, 14).Value = -108180  # N110: -75455 -> -108180

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 248.75  # H19: 334573.47 -> 248.75
$ws.Cells.Item(19, 9).Value = 40.833332  # I19: 358400.5 -> 40.833332
$ws.Cells.Item(19, 10).Value = 872.5  # J19: 995 -> 872.5
$ws.Cells.Item(19, 11).Value = 40.833332  # K19: 358400.5 -> 40.833332
$ws.Cells.Item(19, 12).Value = 872.5  # L19: 995 -> 872.5
$ws.Cells.Item(19, 13).Value = 129.166668  # M19: -358230.5 -> 129.166668
$ws.Cells.Item(19, 14).Value = -1212.5  # N19: -1335 -> -1212.5

$ws.Cells.Item(24, 8).Value = 248.75  # H24: 334573.47 -> 248.75
$ws.Cells.Item(24, 9).Value = 40.833332  # I24: 358400.5 -> 40.833332
$ws.Cells.Item(24, 10).Value = 872.5  # J24: 995 -> 872.5
$ws.Cells.Item(24, 11).Value = 40.833332  # K24: 358400.5 -> 40.833332
$ws.Cells.Item(24, 12).Value = 872.5  # L24: 995 -> 872.5
$ws.Cells.Item(24, 13).Value = 129.166668  # M24: -358230.5 -> 129.166668
$ws.Cells.Item(24, 14).Value = -1212.5  # N24: -1335 -> -1212.5

$ws.Cells.Item(31, 8).Value = 12616.583  # H31: 12408 -> 12616.583
$ws.Cells.Item(31, 9).Value = 5714.9585  # I31: 5582.2 -> 5714.9585
$ws.Cells.Item(31, 11).Value = 5714.9585  # K31: 5582.2 -> 5714.9585
$ws.Cells.Item(31, 13).Value = -5419.9585  # M31: -5287.2 -> -5419.9585

$ws.Cells.Item(34, 8).Value = 12616.583  # H34: 12408 -> 12616.583
$ws.Cells.Item(34, 9).Value = 5714.9585  # I34: 5582.2 -> 5714.9585
$ws.Cells.Item(34, 11).Value = 5714.9585  # K34: 5582.2 -> 5714.9585
$ws.Cells.Item(34, 13).Value = -5512.9585  # M34: -5380.2 -> -5512.9585

$ws.Cells.Item(36, 8).Value = 500  # H36: 0 -> 500
$ws.Cells.Item(36, 9).Value = 500  # I36: 0 -> 500
$ws.Cells.Item(36, 11).Value = 500  # K36: 0 -> 500
$ws.Cells.Item(36, 13).Value = -112  # M36: (new cell) -> -112

$ws.Cells.Item(40, 8).Value = 500  # H40: 0 -> 500
$ws.Cells.Item(40, 9).Value = 500  # I40: 0 -> 500
$ws.Cells.Item(40, 11).Value = 500  # K40: 0 -> 500
$ws.Cells.Item(40, 13).Value = -340  # M40: (new cell) -> -340

$ws.Cells.Item(132, 8).Value = 6546.9644  # H132: 6765.3706 -> 6546.9644
$ws.Cells.Item(132, 9).Value = 1962.3889  # I132: 2039.5883 -> 1962.3889
$ws.Cells.Item(132, 11).Value = 5887.1667  # K132: 6118.7649 -> 5887.1667
$ws.Cells.Item(132, 13).Value = -3357.1667  # M132: -3588.7649 -> -3357.1667

$ws.Cells.Item(134, 8).Value = 52641630  # H134: 41674850 -> 52641630
$ws.Cells.Item(134, 9).Value = 2406.5715  # I134: 1942.8334 -> 2406.5715
$ws.Cells.Item(134, 10).Value = 83347850  # J134: 83347750 -> 83347850
$ws.Cells.Item(134, 11).Value = 7219.7145  # K134: 5828.5002 -> 7219.7145
$ws.Cells.Item(134, 12).Value = 250043550  # L134: 250043250 -> 250043550
$ws.Cells.Item(134, 13).Value = -4684.7145  # M134: -3293.5002 -> -4684.7145
$ws.Cells.Item(134, 14).Value = -250048620  # N134: -250048320 -> -250048620

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 57064256  # H4: 46688960 -> 57064256
$ws.Cells.Item(4, 10).Value = 1000000  # J4: 200092.4 -> 1000000
$ws.Cells.Item(4, 12).Value = 3000000  # L4: 600277.2 -> 3000000
$ws.Cells.Item(4, 14).Value = -3000224  # N4: -600501.2 -> -3000224

$ws.Cells.Item(131, 8).Value = 1491.42  # H131: 1495.52 -> 1491.42
$ws.Cells.Item(131, 10).Value = 1493.9485  # J131: 1498.1753 -> 1493.9485
$ws.Cells.Item(131, 12).Value = 4481.845499999999  # L131: 4494.525900000001 -> 4481.845499999999
$ws.Cells.Item(131, 14).Value = -14561.8455  # N131: -14574.5259 -> -14561.8455

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 836.5  # H2: 837.4643 -> 836.5
$ws.Cells.Item(2, 9).Value = 1146.8422  # I2: 1148.2632 -> 1146.8422
$ws.Cells.Item(2, 11).Value = 1146.8422  # K2: 1148.2632 -> 1146.8422
$ws.Cells.Item(2, 13).Value = -1033.8422  # M2: -1035.2632 -> -1033.8422

$ws.Cells.Item(102, 8).Value = 3887.48  # H102: 3927.92 -> 3887.48
$ws.Cells.Item(102, 9).Value = 4119.864  # I102: 4165.8184 -> 4119.864
$ws.Cells.Item(102, 11).Value = 4119.864  # K102: 4165.8184 -> 4119.864
$ws.Cells.Item(102, 13).Value = -2497.864  # M102: -2543.8184 -> -2497.864

$ws.Cells.Item(135, 8).Value = 154081.33  # H135: 152915.44 -> 154081.33
$ws.Cells.Item(135, 10).Value = 154081.33  # J135: 152915.44 -> 154081.33
$ws.Cells.Item(135, 12).Value = 154081.33  # L135: 152915.44 -> 154081.33
$ws.Cells.Item(135, 14).Value = -164221.33  # N135: -163055.44 -> -164221.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 11759.5  # H7: 11777.1 -> 11759.5
$ws.Cells.Item(7, 9).Value = 15119  # I7: 16032.444 -> 15119
$ws.Cells.Item(7, 10).Value = 8400  # J7: 8295.454 -> 8400
$ws.Cells.Item(7, 11).Value = 15119  # K7: 16032.444 -> 15119
$ws.Cells.Item(7, 12).Value = 8400  # L7: 8295.454 -> 8400
$ws.Cells.Item(7, 13).Value = -15007  # M7: -15920.444 -> -15007
$ws.Cells.Item(7, 14).Value = -8624  # N7: -8519.454 -> -8624

$ws.Cells.Item(22, 8).Value = 4523.3335  # H22: 4363.1816 -> 4523.3335
$ws.Cells.Item(22, 10).Value = 12561.6  # J22: 10634.667 -> 12561.6
$ws.Cells.Item(22, 12).Value = 12561.6  # L22: 10634.667 -> 12561.6
$ws.Cells.Item(22, 14).Value = -13151.6  # N22: -11224.667 -> -13151.6

$ws.Cells.Item(27, 8).Value = 4523.3335  # H27: 4363.1816 -> 4523.3335
$ws.Cells.Item(27, 10).Value = 12561.6  # J27: 10634.667 -> 12561.6
$ws.Cells.Item(27, 12).Value = 12561.6  # L27: 10634.667 -> 12561.6
$ws.Cells.Item(27, 14).Value = -12775.6  # N27: -10848.667 -> -12775.6

$ws.Cells.Item(126, 8).Value = 11759.5  # H126: 11777.1 -> 11759.5
$ws.Cells.Item(126, 9).Value = 15119  # I126: 16032.444 -> 15119
$ws.Cells.Item(126, 10).Value = 8400  # J126: 8295.454 -> 8400
$ws.Cells.Item(126, 11).Value = 45357  # K126: 48097.33199999999 -> 45357
$ws.Cells.Item(126, 12).Value = 25200  # L126: 24886.362 -> 25200
$ws.Cells.Item(126, 13).Value = -42887  # M126: -45627.33199999999 -> -42887
$ws.Cells.Item(126, 14).Value = -30140  # N126: -29826.362 -> -30140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 6158.0728  # H132: 5024.5557 -> 6158.0728
$ws.Cells.Item(132, 9).Value = 3526.1316  # I132: 2855.7637 -> 3526.1316
$ws.Cells.Item(132, 11).Value = 10578.3948  # K132: 8567.2911 -> 10578.3948
$ws.Cells.Item(132, 13).Value = -8048.3948  # M132: -6037.2911 -> -8048.3948
